$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Status" column (G). Build header + first data row (with the new
#    green-fill style) first, since that is the order the style table (and
#    the new cellXfs/fill entries) need to appear in.
# ---------------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Status"

$ws.Range("G2").Value = "Passed"
$ws.Range("G2").Interior.ColorIndex = 17

# Propagate the new "Passed" style down the rest of the Status column.
$ws.Range("G2").Copy()
$ws.Range("G3:G13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Propagate column A:F formatting from row 2 down to the new rows 3:13.
# ---------------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the straightforward columns: InitialDeposit (A), TermMonths (B),
#    Compounding (D), ExpectedValue (E), Status (G, value already "Passed").
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 50000
$ws.Range("B3").Value = 12
$ws.Range("D3").Value = "Compounded Daily"
$ws.Range("E3").Value = 51010.04

$ws.Range("A4").Value = 25000
$ws.Range("B4").Value = 24
$ws.Range("D4").Value = "Compounded Daily"
$ws.Range("E4").Value = 25916.37

$ws.Range("A5").Value = 100000
$ws.Range("B5").Value = 60
$ws.Range("D5").Value = "Compounded Daily"
$ws.Range("E5").Value = 113314.36

$ws.Range("A6").Value = 15000
$ws.Range("B6").Value = 36
$ws.Range("D6").Value = "Compounded Daily"
$ws.Range("E6").Value = 15643.4

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = "Compounded Daily"
$ws.Range("E7").Value = 1

$ws.Range("A8").Value = 9999999
$ws.Range("B8").Value = 120
$ws.Range("D8").Value = "Compounded Daily"
$ws.Range("E8").Value = 16486646.49

$ws.Range("A9").Value = 30000
$ws.Range("B9").Value = 0
$ws.Range("D9").Value = "Compounded Daily"
$ws.Range("E9").Value = 16486646.49

$ws.Range("A10").Value = 30000
$ws.Range("B10").Value = -12
$ws.Range("D10").Value = "Compounded Daily"
$ws.Range("E10").Value = 30483.85

$ws.Range("A11").Value = 30000
$ws.Range("B11").Value = 12
$ws.Range("D11").Value = "Compounded Daily"
$ws.Range("E11").Value = 30483.85

$ws.Range("A12").Value = 30000
$ws.Range("B12").Value = 12
$ws.Range("D12").Value = "Compounded Daily"
$ws.Range("E12").Value = 30301.5

$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 12
$ws.Range("D13").Value = "Compounded Daily"
$ws.Range("E13").Value = 30301.5

# ---------------------------------------------------------------------------
# 4) Notes column (F) - each row gets a brand-new shared string, written in
#    row order (rows 3 through 13).
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "Short-term CD"
$ws.Range("F4").Value = "Medium term"
$ws.Range("F5").Value = "Large principal"
$ws.Range("F6").Value = "Low deposit"
$ws.Range("F7").Value = "Minimum edge case"
$ws.Range("F8").Value = "Max edge case"
$ws.Range("F9").Value = "Zero months"
$ws.Range("F10").Value = "Negative term"
$ws.Range("F11").Value = "Zero interest"
$ws.Range("F12").Value = "Negative interest"
$ws.Range("F13").Value = "Zero deposit"

# ---------------------------------------------------------------------------
# 5) InterestRate column (C) - new unique text values, in first-use order.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "02.00"
$ws.Range("C4").Value = "01.80"
$ws.Range("C5").Value = "02.50"
$ws.Range("C6").Value = "01.40"
$ws.Range("C8").Value = "05.00"
$ws.Range("C12").Value = "-01.00"
$ws.Range("C7").Value = "00.01"

# Reuse already-existing shared strings.
$ws.Range("C9").Value = "01.60"
$ws.Range("C10").Value = "01.60"
$ws.Range("C13").Value = "02.00"

# Zero-interest edge case is stored as a literal number, not text.
$ws.Range("C11").Value = 0

Write-Host "edit complete"
